$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Remove the "Picture 2" picture shape (achievement1_square.jpg) that is no
# longer needed, per the fix for the goals/dashboard pages.
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Name -eq "Picture 2") {
        $shape.Delete()
    }
}
